$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------------
# Change 1: update the Figma prototype link URL (node-id + t query params)
# ---------------------------------------------------------------------------
$oldUrl = "https://www.figma.com/proto/NKvmhTcLxv9uYGuXchTnCH/Desktop-TechWebsite?node-id=2125-343&t=HrFnzWZK1fXnm1g6-1"
$newUrl = "https://www.figma.com/proto/NKvmhTcLxv9uYGuXchTnCH/Desktop-TechWebsite?node-id=2186-73&t=U7VqhBZ2sM35UjpE-1"
$d.Content.Find.Execute($oldUrl, $true, $false, $false, $false, $false, $true, 1, $false, $newUrl, 2) | Out-Null

# Remove the now-stray empty paragraph that used to follow the link paragraph
$linkPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd("`r","`a") -eq $newUrl) { $linkPara = $p; break }
}
$emptyAfterLink = $linkPara.Next()
if ($emptyAfterLink.Range.Text.TrimEnd("`r","`a") -eq "") {
    $emptyAfterLink.Range.Delete() | Out-Null
}

# ---------------------------------------------------------------------------
# Change 2: move <w:lastRenderedPageBreak/> from the "BUY BUTTON" run to the
# following "Font-" run (first BUY BUTTON/Font- pair, under SLIDER LEARN MORE)
# ---------------------------------------------------------------------------
$buyButtonPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd("`r","`a") -eq "BUY BUTTON:  BG- White (C 35)") { $buyButtonPara = $p; break }
}
$fontPara = $buyButtonPara.Next()

$pairRange = $d.Range($buyButtonPara.Range.Start, $fontPara.Range.End)
$pairXml = @"
<w:p $wns><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>BUY BUTTON</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">:  BG- </w:t></w:r><w:r><w:t>White (C 35)</w:t></w:r></w:p><w:p $wns><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:lastRenderedPageBreak/><w:t xml:space="preserve">                                 Font- </w:t></w:r><w:r><w:t>Calibri 40 (B), 003185</w:t></w:r></w:p>
"@
$pairRange.InsertXML($pairXml)

# ---------------------------------------------------------------------------
# Change 3: add the new "QR CODE AFTER CHECKOUT" block at the end of the doc,
# right after "TEXT FONT- Poppins 16 (B)" and before the trailing blank paras
# ---------------------------------------------------------------------------
$anchorPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd("`r","`a") -eq "TEXT FONT- Poppins 16 (B)") { $anchorPara = $p }
}
# Create a fresh, genuinely-empty paragraph right after the anchor first -
# inserting XML straight onto a collapsed range tends to consume/replace the
# paragraph that follows instead of splicing in a new one, so we materialize
# the target paragraph explicitly before swapping its (empty) content out.
$anchorPara.Range.InsertParagraphAfter() | Out-Null
$insertTarget = $anchorPara.Next()
$blockXml = @"
<w:p $wns/>
<w:p $wns>
  <w:pPr>
    <w:rPr>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="50"/>
      <w:szCs w:val="50"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="50"/>
      <w:szCs w:val="50"/>
    </w:rPr>
    <w:lastRenderedPageBreak/>
    <w:t>QR CODE AFTER CHECKOUT</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="50"/>
      <w:szCs w:val="50"/>
    </w:rPr>
    <w:t>:</w:t>
  </w:r>
</w:p>
<w:p $wns>
  <w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>BG-</w:t></w:r>
  <w:r><w:t xml:space="preserve"> </w:t></w:r>
  <w:r><w:t>412</w:t></w:r>
  <w:r><w:t xml:space="preserve"> X </w:t></w:r>
  <w:r><w:t>500</w:t></w:r>
  <w:r><w:t>, White</w:t></w:r>
</w:p>
<w:p $wns>
  <w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>QR-</w:t></w:r>
  <w:r><w:t xml:space="preserve"> YOUR CHOICE</w:t></w:r>
</w:p>
<w:p $wns>
  <w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">CANCEL </w:t></w:r>
  <w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">BUTTON- </w:t></w:r>
  <w:r><w:t>304 X 50</w:t></w:r>
  <w:r><w:t xml:space="preserve">, Any </w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:t>Color</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
</w:p>
<w:p $wns/>
<w:p $wns/>
"@
$insertTarget.Range.InsertXML($blockXml)
